# Trade #79 closed at 2026-02-17 15:49:43 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade #79 (MarketMaking strategy).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - roll-up metrics
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.96   # Current Capital
$wsSummary.Range("B4").Value = -0.05     # Total P&L $
$wsSummary.Range("B5").Value = -0.01     # Total P&L %
$wsSummary.Range("B6").Value = 79        # Total Trades
$wsSummary.Range("B7").Value = 26        # Winning Trades
$wsSummary.Range("B9").Value = 32.91     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.95999999999999   # Capital
$wsStatus.Range("D4").Value = 79                  # Trades
$wsStatus.Range("E4").Value = -0.05               # P&L $
$wsStatus.Range("F4").Value = -0.04               # P&L %
$wsStatus.Range("G4").Value = 32.91               # Win Rate %

# ---------------------------------------------------------------------
# 3. Helper to append the new trade row (#79) to a trade log sheet
# ---------------------------------------------------------------------
function Add-TradeRow($ws, $row) {
    $ws.Range("A$row").Value = 79

    # Keep the date/time text literal - otherwise Excel's type inference
    # would turn "2026-02-17" into a date serial number.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = "2026-02-17"
    $ws.Range("B$row").ClearFormats()

    $ws.Range("C$row").Value = "15:49:37"

    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 0.92
    $ws.Range("G$row").Value = 0.98
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = 6.5217
    $ws.Range("J$row").Value = 0.06
    $ws.Range("K$row").Value = 99.95999999999999
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.15
}

# ---------------------------------------------------------------------
# 4. All Trades sheet - append row 80
# ---------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades 80

# ---------------------------------------------------------------------
# 5. MarketMaking sheet - append row 80
# ---------------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking 80
